# Apply the "Updated symbol list" data refresh to the crypto price table.
# Rows 8-17 also shift: the coin roster inserted GateToken at the top of that
# block, pushing MXToken / LiechtensteinCryptoassetsExchange / WazirX / ... down
# by one row each, so B/C (name/link) are rewritten there too, not just D/E.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'291.03"
$ws.Range("E2").Value = "'0.21%"

# Row 3
$ws.Range("D3").Value = "'30.90"
$ws.Range("E3").Value = "'0.46%"

# Row 4
$ws.Range("D4").Value = "'4.941"
$ws.Range("E4").Value = "'1.41%"

# Row 5
$ws.Range("D5").Value = "'0.07427"
$ws.Range("E5").Value = "'2.54%"

# Row 6
$ws.Range("D6").Value = "'2.195"
$ws.Range("E6").Value = "'-6.55%"

# Row 7
$ws.Range("D7").Value = "'7.713"
$ws.Range("E7").Value = "'0.78%"

# Row 8: GateToken
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'3.750"
$ws.Range("E8").Value = "'1.42%"

# Row 9: MXToken
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9185"
$ws.Range("E9").Value = "'2.30%"

# Row 10: LiechtensteinCryptoassetsExchange
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.09238"
$ws.Range("E10").Value = "'14.67%"

# Row 11: WazirX
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1706"
$ws.Range("E11").Value = "'1.95%"

# Row 12: MandalaExchangeToken
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08277"
$ws.Range("E12").Value = "'1.53%"

# Row 13: BitrueCoin
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03218"
$ws.Range("E13").Value = "'4.73%"

# Row 14: BitMartToken
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09992"
$ws.Range("E14").Value = "'-0.33%"

# Row 15: BitForexToken
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001497"
$ws.Range("E15").Value = "'0.15%"

# Row 16: TigerCash
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005805"
$ws.Range("E16").Value = "'0.54%"

# Row 17: LEO
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.478"
$ws.Range("E17").Value = "'0.01%"

# Row 18
$ws.Range("D18").Value = "'2.060"
$ws.Range("E18").Value = "'-0.83%"

# Row 19
$ws.Range("D19").Value = "'0.3330"
$ws.Range("E19").Value = "'0.47%"

# Row 20
$ws.Range("D20").Value = "'0.1288"
$ws.Range("E20").Value = "'0.04%"

# Row 21
$ws.Range("D21").Value = "'4.155"
$ws.Range("E21").Value = "'4.64%"

# Row 22
$ws.Range("D22").Value = "'0.2120"
$ws.Range("E22").Value = "'0.55%"

# Row 23
$ws.Range("D23").Value = "'0.04507"
$ws.Range("E23").Value = "'-0.25%"

# Row 24
$ws.Range("D24").Value = "'0.001216"
$ws.Range("E24").Value = "'0.16%"

# Row 25
$ws.Range("D25").Value = "'0.004235"
$ws.Range("E25").Value = "'-4.00%"

# Row 26
$ws.Range("D26").Value = "'0.0001298"
$ws.Range("E26").Value = "'-0.16%"

# Row 27
$ws.Range("D27").Value = "'0.0003388"
$ws.Range("E27").Value = "'-0.25%"

# Row 39
$ws.Range("D39").Value = "'0.01590"
$ws.Range("E39").Value = "'0.16%"

# Row 40
$ws.Range("D40").Value = "'0.04562"
$ws.Range("E40").Value = "'4.20%"

# Row 41
$ws.Range("D41").Value = "'0.007374"
$ws.Range("E41").Value = "'1.09%"

# Row 42
$ws.Range("D42").Value = "'0.009846"
$ws.Range("E42").Value = "'-1.75%"

# Row 43
$ws.Range("D43").Value = "'0.1345"
$ws.Range("E43").Value = "'2.42%"

# Row 44
$ws.Range("D44").Value = "'0.002157"
$ws.Range("E44").Value = "'3.68%"

# Row 45
$ws.Range("D45").Value = "'0.008529"
$ws.Range("E45").Value = "'-7.11%"

# Row 46
$ws.Range("D46").Value = "'0.00006117"
$ws.Range("E46").Value = "'7.25%"

# Row 47
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.25%"

# Row 48
$ws.Range("D48").Value = "'2.594"
$ws.Range("E48").Value = "'15.68%"

# Row 49
$ws.Range("D49").Value = "'0.001996"
$ws.Range("E49").Value = "'-31.17%"

# Row 50
$ws.Range("D50").Value = "'0.00002096"
$ws.Range("E50").Value = "'-0.25%"

# Row 51
$ws.Range("D51").Value = "'0.0001996"
$ws.Range("E51").Value = "'-0.25%"
